$d = $word.ActiveDocument

$enDash = [char]0x2013
$oldTitleText = "Play Age of the Gods " + $enDash + " Fate Sister for Free | Exciting Bonus Features"
$oldMetaDescText = "Find out why Age of the Gods " + $enDash + " Fate Sister is worth playing! Enjoy exciting bonus features and win progressive jackpots. Play for free!"
$metaDescSuffix = ": " + $oldMetaDescText
$newImagePromptText = "Create a feature image that captures the mythical world of Age of the Gods " + $enDash + " Fate Sisters, featuring a happy Maya warrior with glasses in cartoon-style. The background should feature an ancient temple with the Fate Sisters standing in a line, with Clotho at the front holding a spindle and Lachesis holding a pair of scissors. Atropos should be shown holding a glowing orb. The Maya warrior should be positioned in front of the Fate Sisters with a big grin on his face, wearing a feathered headdress and glasses. The overall tone of the image should be bright and colorful to capture the excitement and adventure of the game."

# ---------------------------------------------------------------------------
# Locate the duplicate bold "Play Age of the Gods..." paragraph near the end
# of the document (it mirrors the H1 title, just bold and un-styled). We
# copy *it* - rather than typing fresh text - so the new paragraph we create
# reuses its exact run layout (a leading empty run followed by the bold run).
# ---------------------------------------------------------------------------
$dupTitleIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq ($oldTitleText + "`r")) {
        $dupTitleIdx = $i
    }
}

$srcPara = $d.Paragraphs($dupTitleIdx)
$srcPara.Range.Copy()

# ---------------------------------------------------------------------------
# Insert a brand-new (Normal-styled) paragraph right after the H1 heading,
# then paste the copied paragraph into it.
# ---------------------------------------------------------------------------
$h1 = $d.Paragraphs(1)
$h1.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$insertPoint = $metaPara.Range
$insertPoint.Collapse(1)
$insertPoint.Paste()

# ---------------------------------------------------------------------------
# Turn the pasted "Play Age of the Gods..." bold run into "Meta description"
# and append the (non-bold) description text after it.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Find.Execute($oldTitleText, $true, $false, $false, $false, $false, $true, 1, $false, "Meta description", 2)

$metaPara = $d.Paragraphs(2)
$metaRange = $metaPara.Range
$appendRange = $d.Range($metaRange.Start, $metaRange.End - 1)
$appendRange.InsertAfter($metaDescSuffix)

# ---------------------------------------------------------------------------
# Remove the now-redundant original duplicate bold title paragraph further
# down the document (its index shifted by +1 after the insertion above).
# ---------------------------------------------------------------------------
$dupTitleIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($i -ne 2 -and $d.Paragraphs($i).Range.Text -eq ($oldTitleText + "`r")) {
        $dupTitleIdx = $i
    }
}
$d.Paragraphs($dupTitleIdx).Range.Delete()

# ---------------------------------------------------------------------------
# Replace the final italic paragraph's text (the old meta-description line)
# with the new AI image-prompt text, keeping its italic formatting.
# ---------------------------------------------------------------------------
$italicIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq ($oldMetaDescText + "`r")) {
        $italicIdx = $i
    }
}
$italicPara = $d.Paragraphs($italicIdx)
$italicPara.Range.Find.Execute($oldMetaDescText, $true, $false, $false, $false, $false, $true, 1, $false, $newImagePromptText, 2)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
